# PersonnelFileRequestLetterDATemplate.docx -- "Implemented Oct 1 feedback"
#
# 1. Merge the "work" and "on_leave" branches of the employment_status
#    conditional: the {% elif employment_status == "work" %} clause grows
#    an "or employment_status == 'on_leave'" condition.
# 2. The (now redundant) dedicated "on_leave" branch is dropped, and the
#    "laid_off" branch is reworded / now references {{end_date}} (the date
#    the person was laid off) instead of repeating {{hiring_date}}.
# 3. A stale <w:lastRenderedPageBreak/> hint (left over from a previous
#    render pass) is cleared from the "{%if add_signature..." run.
#
# The Jinja-ish template text uses curly ("smart") quotes around the
# work/on_leave/laid_off string literals, but plain straight quotes inside
# employer.name.full(middle='full'). We therefore assign directly to
# Range.Text (rather than going through Find.Replacement.Text) so Word's
# AutoCorrect "smart quotes" feature does not also convert those straight
# quotes.

$q1 = [char]0x201C   # “
$q2 = [char]0x201D   # ”

$d = $word.ActiveDocument

function Replace-FirstMatch([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.MatchWildcards = $false
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $start = $rng.Start
    $end = $rng.End
    $target = $d.Range($start, $end)
    $target.Text = $newText
}

# --- Change 1: "work" elif also covers "on_leave" -----------------------
$oldWork = $q1 + "work" + $q2 + " %}I have been"
$newWork = $q1 + "work" + $q2 + " or employment_status == " + $q1 + "on_leave" + $q2 + " %}I have been"
Replace-FirstMatch $oldWork $newWork

# --- Change 2: rewrite the "on_leave" / "laid_off" branches --------------
$oldTail = "{% elif employment_status == " + $q1 + "on_leave" + $q2 + " %} I have been an employee of {{ employer.name.full(middle='full')}} since around {{hiring_date}} and am currently on a leave of absence. {% elif employment_status == " + $q1 + "laid_off" + $q2 + " %} I have been an employee of {{ employer.name.full(middle='full')}} since around {{hiring_date}}. I was laid off, but I am subject to recall."
$newTail = "{% elif employment_status == " + $q1 + "laid_off" + $q2 + " %} I started working as an employee of {{ employer.name.full(middle='full')}} around {{hiring_date}}. I was laid off around {{end_date}}, but I am subject to recall."
Replace-FirstMatch $oldTail $newTail

# --- Change 3: drop the stale lastRenderedPageBreak hint before "{%if add"
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found3 = $rng3.Find.Execute("{%if add", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $s3 = $rng3.Start
    $e3 = $rng3.End
    $tmp = $d.Range($s3, $e3)
    $tmp.Text = "{%if ZZZADD_TMP"
    $s3b = $s3
    $e3b = $s3 + 15
    $tmp2 = $d.Range($s3b, $e3b)
    $tmp2.Text = "{%if add"
}

Write-Output "Done"
